# Updated symbol list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values (price, volume %, hour) must stay stored as TEXT,
# matching the source data's inline-string cell type. Assigning via .Value
# on a numeric-looking string would auto-coerce it to a Number, so we write
# it as a quote-prefixed Formula (forces Text) and then reset the cell style
# back to Normal so no stray number-format/style is left behind.
$cellText = @{
    "D2" = '315.69'
    "E2" = '3.19%'
    "G2" = '13'
    "D3" = '39.52'
    "E3" = '2.87%'
    "G3" = '13'
    "D4" = '5.133'
    "E4" = '0.72%'
    "G4" = '13'
    "D5" = '0.08210'
    "E5" = '1.72%'
    "G5" = '13'
    "E6" = '0.99%'
    "G6" = '13'
    "E7" = '3.60%'
    "G7" = '13'
    "D8" = '0.9294'
    "E8" = '0.10%'
    "G8" = '13'
    "D9" = '0.1400'
    "E9" = '-3.50%'
    "G9" = '13'
    "D10" = '0.1982'
    "E10" = '2.43%'
    "G10" = '13'
    "D11" = '0.09106'
    "E11" = '0.59%'
    "G11" = '13'
    "D12" = '0.03509'
    "E12" = '-0.14%'
    "G12" = '13'
    "D13" = '0.09815'
    "E13" = '0.15%'
    "G13" = '13'
    "D14" = '0.001394'
    "E14" = '-1.56%'
    "G14" = '13'
    "D15" = '0.006101'
    "E15" = '0.13%'
    "G15" = '13'
    "D16" = '3.657'
    "E16" = '-1.81%'
    "G16" = '13'
    "D17" = '4.261'
    "E17" = '1.17%'
    "G17" = '13'
    "D18" = '3.298'
    "E18" = '-3.45%'
    "G18" = '13'
    "D19" = '0.3463'
    "E19" = '0.06%'
    "G19" = '13'
    "D20" = '0.1294'
    "E20" = '-3.00%'
    "G20" = '13'
    "D21" = '4.921'
    "E21" = '2.10%'
    "G21" = '13'
    "D22" = '0.2447'
    "E22" = '-0.44%'
    "G22" = '13'
    "D23" = '0.04326'
    "E23" = '-1.10%'
    "G23" = '13'
    "E24" = '-0.65%'
    "G24" = '13'
    "D25" = '0.004785'
    "E25" = '15.89%'
    "G25" = '13'
    "D26" = '0.0001296'
    "E26" = '-0.66%'
    "G26" = '13'
    "D27" = '0.0003994'
    "E27" = '-10.19%'
    "G27" = '13'
    "G28" = '13'
    "G29" = '13'
    "G30" = '13'
    "G31" = '13'
    "G32" = '13'
    "G33" = '13'
    "G34" = '13'
    "G35" = '13'
    "G36" = '13'
    "G37" = '13'
    "G38" = '13'
    "D39" = '0.02225'
    "E39" = '7.52%'
    "G39" = '13'
    "D40" = '0.05266'
    "E40" = '4.47%'
    "G40" = '13'
    "D41" = '0.007556'
    "E41" = '1.21%'
    "G41" = '13'
    "D42" = '0.009788'
    "E42" = '-3.22%'
    "G42" = '13'
    "D43" = '0.1376'
    "E43" = '1.74%'
    "G43" = '13'
    "D44" = '0.002123'
    "E44" = '-1.13%'
    "G44" = '13'
    "D45" = '0.009808'
    "E45" = '6.51%'
    "G45" = '13'
    "D46" = '0.00006360'
    "E46" = '2.77%'
    "G46" = '13'
    "D47" = '0.00000000749'
    "E47" = '-0.39%'
    "G47" = '13'
    "D48" = '0.002764'
    "E48" = '-7.62%'
    "G48" = '13'
    "D49" = '0.001198'
    "E49" = '-25.26%'
    "G49" = '13'
    "D50" = '0.00002097'
    "E50" = '-0.39%'
    "G50" = '13'
    "D51" = '0.0001997'
    "E51" = '-0.39%'
    "G51" = '13'
}
foreach ($ref in $cellText.Keys) {
    $cell = $ws.Range($ref)
    $cell.Formula = "'" + $cellText[$ref]
    $cell.Style = "Normal"
}

# Plain text values (coin names / links) are already non-numeric, so a
# normal .Value assignment already stores them as Text with no style change.
$cellPlain = @{
    "B48" = "BOLO"
    "C48" = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
    "B49" = "CoinbaseStockToken"
    "C49" = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
}
foreach ($ref in $cellPlain.Keys) {
    $ws.Range($ref).Value = $cellPlain[$ref]
}
